$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.300.05'
$ws.Range('E2').Value = '  +0.31%  '
$ws.Range('D3').Value = '2.008.79'
$ws.Range('E3').Value = '  -0.07%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '258.29'
$ws.Range('E5').Value = '  +4.59%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  -1.61%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '55.94'
$ws.Range('E8').Value = '  -6.38%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.385'
$ws.Range('E9').Value = '  -0.79%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0768'
$ws.Range('E10').Value = '  -5.42%  '
$ws.Range('E11').Value = '  -2.01%  '
$ws.Range('D12').Value = '2.302.47'
$ws.Range('E12').Value = '  -0.01%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '14.21'
$ws.Range('E13').Value = '  -5.79%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '21.04'
$ws.Range('E14').Value = '  -5.87%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.800'
$ws.Range('E15').Value = '  -5.41%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '5.23'
$ws.Range('E16').Value = '  -4.49%  '
$ws.Range('D17').Value = '2.009.68'
$ws.Range('E17').Value = '  -0.05%  '
$ws.Range('D18').Value = '37.197.23'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '69.62'
$ws.Range('E20').Value = '  -3.76%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '5.12'
$ws.Range('E21').Value = '  -2.01%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '228.29'
$ws.Range('E22').Value = '  -1.08%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.61'
$ws.Range('E23').Value = '  +4.86%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '1.00'
$ws.Range('E24').Value = '  +0.02%  '
$ws.Range('E25').Value = '  -0.51%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.71'
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  -6.29%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '19.60'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.130'
$ws.Range('E29').Value = '  -7.24%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '1.32'
$ws.Range('E30').Value = '  -4.05%  '
$ws.Range('E32').Value = '  -4.10%  '
$ws.Range('E33').Value = '  -2.31%  '
$ws.Range('E34').Value = '  +1.33%  '
$ws.Range('E35').Value = '  -3.68%  '
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('E37').Value = '  +0.08%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.34'
$ws.Range('E38').Value = '  -3.67%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '5.26'
$ws.Range('E39').Value = '  -1.27%  '
$ws.Range('E40').Value = '  +3.71%  '
$ws.Range('E41').Value = '  +0.94%  '
$ws.Range('E42').Value = '  -5.53%  '
$ws.Range('E43').Value = '  -1.34%  '
$ws.Range('D44').Value = '1.397.54'
$ws.Range('E44').Value = '  +1.80%  '
$ws.Range('B45').Value = 'InjectiveProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '15.69'
$ws.Range('E45').Value = '  -5.63%  '
$ws.Range('B46').Value = 'Aave'
$ws.Range('C46').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '89.41'
$ws.Range('E46').Value = '  -2.97%  '
$ws.Range('E47').Value = '  -2.85%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '7.03'
$ws.Range('E48').Value = '  -5.03%  '
$ws.Range('B49').Value = 'MXToken'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.91'
$ws.Range('E49').Value = '  +2.23%  '
$ws.Range('D50').Value = '2.194.66'
$ws.Range('E50').Value = '  +0.04%  '
$ws.Range('E51').Value = '  -7.25%  '
